# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# F2: 282 -> 283
# F4: 95  -> 96
# F5: 858 -> 860

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 283
    $ws.Range("F4").Value = 96
    $ws.Range("F5").Value = 860
}
